# Update participant data: User_ID (column A) and Rating (column C) values
# for rows 2-16 on the active worksheet, per the source data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  UserId = 3815; Rating = 4 },
    @{ Row = 3;  UserId = 3759; Rating = 2 },
    @{ Row = 4;  UserId = 3247; Rating = 4 },
    @{ Row = 5;  UserId = 8918; Rating = 5 },
    @{ Row = 6;  UserId = 2456; Rating = 3 },
    @{ Row = 7;  UserId = 8013; Rating = 2 },
    @{ Row = 8;  UserId = 8941; Rating = 2 },
    @{ Row = 9;  UserId = 1581; Rating = 4 },
    @{ Row = 10; UserId = 5719; Rating = 3 },
    @{ Row = 11; UserId = 7040; Rating = 2 },
    @{ Row = 12; UserId = 3315; Rating = 4 },
    @{ Row = 13; UserId = 9821; Rating = 4 },
    @{ Row = 14; UserId = 6117; Rating = 3 },
    @{ Row = 15; UserId = 6541; Rating = 5 },
    @{ Row = 16; UserId = 5724; Rating = 4 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.UserId
    $ws.Cells.Item($u.Row, 3).Value = $u.Rating
}
